# Review_321.docx content update
$d = $word.ActiveDocument

# --- Paragraph 1: title line (date + paper title, two runs split by a <w:br/>) ---
$d.Content.Find.Execute(
    'המאמר היומי של מייק -16.10.24:',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'המאמר היומי של מייק -15.10.24:',
    2)

$d.Content.Find.Execute(
    'EFFICIENT REINFORCEMENT LEARNING WITH LARGE LANGUAGE MODEL PRIORS',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'EFFICIENT DICTIONARY LEARNING WITH SWITCH SPARSE AUTOENCODERS',
    2)

# --- Paragraph 2: intro ---
$d.Content.Find.Execute(
    'היום נסקור מאמר שהוא נראה די כבד מתמטית (הרבה נוסחאות ומלל שנראה מתמטי) אבל הרעיון מאחוריו הוא די פשוט וקל להסבר. אנחנו אוהבים למנף  את עוצמתם של מודלי שפה למשימות רבות (ולא תמיד לכאלו שהם מסוגלים לבצע כמו שצריך לפחות כרגע).',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'היום סוקרים מאמר קליל המשלב שני רעיונות די נחמדים שמשמים LLMs (במיוחד לאחרונה) והאמת השילוב שלהם נראה די טבעי. הרעיון הראשון הינו Mixture of Experts או MoE בקצרה.',
    2)

# --- Paragraph 3: MoE explanation ---
$d.Content.Find.Execute(
    'המאמר מציע להשתמש במודל שפה כפריור עבור סוכנים במשימות בהם הם צריכים לבצע SDM או sequential decision making. המאמר נותן בתור דוגמא משחק overcooked כאשר הסוכן צריך לבצע משימות בישול שונות בהתבסס על מצב המטבח שבו הוא מבשל אותם. המטרה של הסוכן היא לחזות את הפעולה הבא (באמצעות תיאור טקסטואלי) כאשר התגמול הוא ביצוע נכון של המשימה (הכנה של מנה לפי המתכון :)).',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'MoE היא שיטה המאפשרת לנו להקל על האינפרנס על ידי שימוש רק בחלק ממשקלי המודל. בד״כ מטריצות משקלים ברשת feed-forward (יש שם 2 שכבות בסך הכל) בבלוק הטרנספורמרים (אחרי attention) מחוקלים לכמה קבוצות שכל אחת מהן נקראת מומחה או expert. באינפרנס המודל משתמש רק בחלק (לפעמים רק אחד) מהמומחים ובכך הוא מוריד את מחירו של האינפרנס. כלומר אותו המודל מופעל בצורה קצת שונה בהתאם לקלט (בנוסף ל-attention),',
    2)

# --- Paragraph 4: SAE intro ---
$d.Content.Find.Execute(
    'כאמור המטרה כאן היא לחזות את הפעולה הבאה עבור הסוכן (המתוארת) על ידי הטקסט כאשר המצב (state) גם מתואר על ידי טקסט. בגדול מאוד אנו מתחילים ממודל אחד (הפריור P) עבור חיזוי המצב הבא (מהמצב הקודם והפעולה) ועבור חיזוי הפעולה הבאה בהינתן המצב (מתואר על ידי התפלגות Q_h). המטרה כאן היא ללמוד את Q_h כאשר ממקסמת התגמול הצפוי ושומרת את התפלגות Q קרובה לפריור P (זוכרים PPO שהתפרסם מאוד לפני שנתיים כאשר OpenAI השתמשו בו ל-RLHF לאימון מודלי שפה). המרחק כמובן ניתן על ידי ה-KL 🙂',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'הקונספט השני הוא Sparse AutoEncoders או SAE בקצרה שהפך להיות די פופולרי אחרי החוקרים של אנטרופיק הציעו להשתמש בו למטרת חקר interpretability של מודלי שפה. לפני הבלוג הזה הסברה הרווחת (סוג של) היתה שבמודל שפה יש נוירונים שנדלקים חזק (מקבלים ערך גבוה) על קונספטים מסוימים כאשר כל נוירון כזה הינו מונו-סמנטי כלומר יש קונספט אחד בלבד שהוא ״אחראי״ עליו. ',
    2)

# --- Paragraph 5: SAE detail ---
$d.Content.Find.Execute(
    'אז הפעולה הבאה a_t (כלומר גנרוט התיאור הטקסטואלי שלה) מתבצע באופן הבא. דוגמים כמה גרסאות של a_t עם P מחשבים את הנראות שלהם לפי Q הנלמד, מנרמלים עם הסופטמקס ודוגמים את הפעולה הבאה כאשר מטרת התהליך מקסום של התגמול הצפוי (עם הרגולריזציה שהסברנו עליה קודם).',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'לעומת זאת החוקרים של אנטרופיק הציע להתבונן בכל נוירון כפולי-סמנטי כלומר ״אחראי״ על מספר קונספטים  לא קשורים. לפי משנתם ניתן לגלות את הקונספטים האלו באמצעות SAE שבונה autoencoder דליל (הרוב אפסים) במימד גבוה הרבה יותר מגודל השכבה שבה נמצאים הנוירונים הפוליסמנטיים אלו. SAE כאן מורכב משתי שכבות בלבד, אחת לאנקודר ואחת לדקודר. ',
    2)

# --- Paragraph 6: monosemantic neuron detail ---
$d.Content.Find.Execute(
    'כמובן שניתן לעשות את זה בכמה אופנים: בצורה של online דרך שערוך של פונקציית Q של הזוג (מצב, פעולה) כאשר פונקציית Q קשורה להתפלגות Q_h של הפעולה הבא שנידונה בפסקה הקודמת (עניין של נרמול נכון). ניתן לעשות את זה גם באמצעות offline עם איזה פוליסי טוב ידוע של המומחים כאשר המטרה היא גם שערוך של פונקציית Q שבאמצעותה ניתן לשערך (לקבל) את Q_h עבור חיזוי הפעולה הבא. ניתן לעשות את זה גם באמצעות שיטה דומה ל-PPO אבל בכל המקרים הפריור הוא ההתפלגות המושרית על ידי מודל שפה נתון.',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ' כאן כל רכיב שהוא לא אפס בווקטור אחרי שכבת ה-encoder של SAE הוא אחראי על קונספט מסוים כלומר מהווה נוירון מונוסמנטי. כך יוצא שכל נוירון בשכבה המקורית הוא שילוב לינארים של הנוירונים המונוסמנטיים אלו. SAE מאומן בצורה די סטנדרטית עם איבר רגולריזציה שאוכף את דלילות הייצוג אחרי האנקודר.',
    2)

# --- Paragraph 7: conclusion ---
$d.Content.Find.Execute(
    'מאמר מעניין בקיצור…',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'אז המאמר מציע לשלב את שני הקונספטים האלו כך שכל נוירון הוא צירוף לינארי אחר של הנוירונים המונוסמנטיים בשכבת ה-encoder. זה מאפשר גמישות נוספת ביחס לרעיון המקורי ובטח מאפשר לגלות קונספטים שונים המוסתרים בתוך ה-LLMs שלנו.',
    2)

# --- Paragraph 8: recommendation line (was the arxiv pdf link) ---
$d.Content.Find.Execute(
    'https://arxiv.org/pdf/2410.07927',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'מאמר קליל - ממליץ להעיף מבט',
    2)

# --- New paragraph 9: the new arxiv abs link, appended after the last paragraph ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'https://arxiv.org/abs/2410.08201'

Write-Output 'edit complete'
